# Update the "parameters" sheet with new IP addresses and the "self" user id,
# then leave "parameters" as the active/selected sheet (matching the
# sheetView/bookView changes captured in the diff).

$wb = $excel.ActiveWorkbook
$wsParams = $wb.Worksheets.Item("parameters")

# Row 3: ip_addr xx.xxx.x.204 -> xx.xxx.x.15, security none -> self
$wsParams.Range("C3").Value = "xx.xxx.x.15"
$wsParams.Range("D3").Value = "self"

# Row 4: ip_addr xx.xxx.x.204 -> xx.xxx.x.26, security none -> self
$wsParams.Range("C4").Value = "xx.xxx.x.26"
$wsParams.Range("D4").Value = "self"

# Make "parameters" the active sheet/tab and set its selection, matching
# the sheetView changes (tabSelected moved from Instructions to parameters).
$wsParams.Select()
$wsParams.Range("F13").Select()
